$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# The roster table (A1:E7) is being reshuffled: "Billy" moves from
# Upstage Right to Upstage Left; the old "Pink/???" placeholder row
# is replaced with a new "Orange / Mike" band member at Upstage
# Right; and Yellow/Blue/Red each shift down one row, picking up new
# stage positions along the way.
# ------------------------------------------------------------------

# --- Shift formatting (fill/font) down the rows before overwriting values,
#     so each row keeps a clean, de-duplicated style where possible ---

# row6 (Red) format -> row7
$ws.Range("A6:E6").Copy() | Out-Null
$ws.Range("A7:E7").PasteSpecial(-4122) | Out-Null

# row5 (Blue) format -> row6
$ws.Range("A5:E5").Copy() | Out-Null
$ws.Range("A6:E6").PasteSpecial(-4122) | Out-Null

# row4 (Yellow) format -> row5
$ws.Range("A4:E4").Copy() | Out-Null
$ws.Range("A5:E5").PasteSpecial(-4122) | Out-Null

# row4 becomes the brand new "Orange" entry: start from a white-font
# template (row2, Green) then recolor the fill to Orange.
$ws.Range("A2:E2").Copy() | Out-Null
$ws.Range("A4:E4").PasteSpecial(-4122) | Out-Null
$ws.Range("A4:E4").Interior.ThemeColor = 6

$excel.CutCopyMode = 0

# --- Cell values ---
$ws.Range("C3").Value = "Upstage Left"
$ws.Range("D3").Value = "USL"

$ws.Range("A4").Value = "Orange"
$ws.Range("B4").Value = "Mike"
$ws.Range("C4").Value = "Upstage Right"
$ws.Range("D4").Value = "USR"
$ws.Range("E4").Value = "Guitar, Vocals"

$ws.Range("A5").Value = "Yellow"
$ws.Range("B5").Value = "Jessica"
$ws.Range("C5").Value = "Downstage Center"
$ws.Range("D5").Value = "DSC"
$ws.Range("E5").Value = "Guitar, Vocals"

$ws.Range("A6").Value = "Blue"
$ws.Range("B6").Value = "Paul"
$ws.Range("C6").Value = "Downstage Right"
$ws.Range("D6").Value = "DSR"
$ws.Range("E6").Value = "Guitar, Vocals"

$ws.Range("A7").Value = "Red"
$ws.Range("B7").Value = "Kenzi"
$ws.Range("C7").Value = "Downstage Left"
$ws.Range("D7").Value = "DSL"
$ws.Range("E7").Value = "Keys, Guitar, Vocals"

# --- Column widths (content got a little narrower in col C, a touch
#     wider in col A, once the new labels took over) ---
$ws.Columns.Item(1).ColumnWidth = 15.9
$ws.Columns.Item(3).ColumnWidth = 39.45

# --- Selection / cursor position left where the editor's cursor ended up ---
$ws.Range("I12").Select() | Out-Null
